$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 09:00:05"
$wsZhCn.Range("H2").Value = "2016-03-21 09:00:29"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 09:00:11"
$wsDeDe.Range("H2").Value = "2016-03-21 09:00:37"
